$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typos/punctuation in Razon social / Nombre Fantasia shared strings (comma -> period) ---
$ws.Cells.Replace("PITTER ROLANDO L.J, CERGNEUX MARIO M. Y CERGNEUX DANIEL F.  S.H.", "PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH", 1) | Out-Null
$ws.Cells.Replace("FERNANDEZ MARIO H, GALLICET OSCAR M", "FERNANDEZ MARIO H. GALLICET OSCAR M", 1) | Out-Null
$ws.Cells.Replace("IZAGUIRRE CARLOS MARIA, MOREND MARIA ELENA Y MOREND MARIA TERESA", "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA", 1) | Out-Null
$ws.Cells.Replace("MARSICO GUILLERMO MIGUEL, MARSICO JUAN EDUARDO", "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO", 1) | Out-Null
$ws.Cells.Replace("RICCOTTI, MARIANA EDITH", "RICCOTTI. MARIANA EDITH", 1) | Out-Null

# --- Fix decimal formatting in Importe column (comma-decimal -> dot-decimal) ---
$importeRange = $ws.Range("H2:H129")
$importeRange.NumberFormat = "@"
$importeRange.Replace("770,00", "770.00", 1) | Out-Null
$importeRange.Replace("2.471,15", "2471.15", 1) | Out-Null
$importeRange.Replace("22.423,54", "22423.54", 1) | Out-Null
$importeRange.Replace("3.223,14", "3223.14", 1) | Out-Null
$importeRange.Replace("410,00", "410.00", 1) | Out-Null
$importeRange.Replace("25.560,00", "25560.00", 1) | Out-Null
$importeRange.Replace("88.796,58", "88796.58", 1) | Out-Null
$importeRange.Replace("12.857,11", "12857.11", 1) | Out-Null
$importeRange.Replace("349,04", "349.04", 1) | Out-Null
$importeRange.Replace("1.796,10", "1796.10", 1) | Out-Null
$importeRange.Replace("2.377,54", "2377.54", 1) | Out-Null
$importeRange.Replace("1.111,00", "1111.00", 1) | Out-Null
$importeRange.Replace("5.347,83", "5347.83", 1) | Out-Null
$importeRange.Replace("80,00", "80.00", 1) | Out-Null
$importeRange.Replace("667,00", "667.00", 1) | Out-Null
$importeRange.Replace("10.636,49", "10636.49", 1) | Out-Null
$importeRange.Replace("1.650,00", "1650.00", 1) | Out-Null
$importeRange.Replace("3.040,00", "3040.00", 1) | Out-Null
$importeRange.Replace("143,00", "143.00", 1) | Out-Null
$importeRange.Replace("332,50", "332.50", 1) | Out-Null
$importeRange.Replace("483,65", "483.65", 1) | Out-Null
$importeRange.Replace("482,13", "482.13", 1) | Out-Null
$importeRange.Replace("18.303,22", "18303.22", 1) | Out-Null
$importeRange.Replace("5.970,99", "5970.99", 1) | Out-Null
$importeRange.Replace("1.413,90", "1413.90", 1) | Out-Null
$importeRange.Replace("74,40", "74.40", 1) | Out-Null
$importeRange.Replace("7.553,00", "7553.00", 1) | Out-Null
$importeRange.Replace("14,08", "14.08", 1) | Out-Null
$importeRange.Replace("799,50", "799.50", 1) | Out-Null
$importeRange.Replace("71,10", "71.10", 1) | Out-Null
$importeRange.Replace("135,00", "135.00", 1) | Out-Null
$importeRange.Replace("800,50", "800.50", 1) | Out-Null
$importeRange.Replace("99,00", "99.00", 1) | Out-Null
$importeRange.Replace("2.139,00", "2139.00", 1) | Out-Null
$importeRange.Replace("966,00", "966.00", 1) | Out-Null
$importeRange.Replace("986,50", "986.50", 1) | Out-Null
$importeRange.Replace("825,00", "825.00", 1) | Out-Null
$importeRange.Replace("374,00", "374.00", 1) | Out-Null
$importeRange.Replace("126,54", "126.54", 1) | Out-Null
$importeRange.Replace("83,60", "83.60", 1) | Out-Null
$importeRange.Replace("1.600,00", "1600.00", 1) | Out-Null
$importeRange.Replace("1.592,00", "1592.00", 1) | Out-Null
$importeRange.Replace("1.495,00", "1495.00", 1) | Out-Null
$importeRange.Replace("459,80", "459.80", 1) | Out-Null
$importeRange.Replace("509,00", "509.00", 1) | Out-Null
$importeRange.Replace("1.700,00", "1700.00", 1) | Out-Null
$importeRange.Replace("54,00", "54.00", 1) | Out-Null
$importeRange.Replace("17.182,80", "17182.80", 1) | Out-Null
$importeRange.Replace("5.540,00", "5540.00", 1) | Out-Null
$importeRange.Replace("2.394,00", "2394.00", 1) | Out-Null
$importeRange.Replace("329.410,00", "329410.00", 1) | Out-Null
$importeRange.Replace("302,00", "302.00", 1) | Out-Null
$importeRange.Replace("29,25", "29.25", 1) | Out-Null
$importeRange.Replace("2.300,08", "2300.08", 1) | Out-Null
$importeRange.Replace("280,00", "280.00", 1) | Out-Null
$importeRange.Replace("1.318,00", "1318.00", 1) | Out-Null
$importeRange.Replace("45,00", "45.00", 1) | Out-Null
$importeRange.Replace("560,00", "560.00", 1) | Out-Null
$importeRange.Replace("680,00", "680.00", 1) | Out-Null
$importeRange.Replace("190,00", "190.00", 1) | Out-Null
$importeRange.Replace("213.732,00", "213732.00", 1) | Out-Null
$importeRange.Replace("4.564,46", "4564.46", 1) | Out-Null
$importeRange.Replace("1.544,63", "1544.63", 1) | Out-Null
$importeRange.Replace("455,73", "455.73", 1) | Out-Null
$importeRange.Replace("1.009,92", "1009.92", 1) | Out-Null
$importeRange.Replace("19.919,18", "19919.18", 1) | Out-Null
$importeRange.Replace("12.119,38", "12119.38", 1) | Out-Null
$importeRange.Replace("510,20", "510.20", 1) | Out-Null
$importeRange.Replace("27,00", "27.00", 1) | Out-Null
$importeRange.Replace("3.601,50", "3601.50", 1) | Out-Null
$importeRange.Replace("4.481,43", "4481.43", 1) | Out-Null
$importeRange.Replace("182,65", "182.65", 1) | Out-Null
$importeRange.Replace("99,64", "99.64", 1) | Out-Null
$importeRange.Replace("10.000,00", "10000.00", 1) | Out-Null
$importeRange.Replace("520,00", "520.00", 1) | Out-Null
$importeRange.Replace("1.238,00", "1238.00", 1) | Out-Null
$importeRange.Replace("1.200,00", "1200.00", 1) | Out-Null
$importeRange.Replace("4.132,23", "4132.23", 1) | Out-Null
$importeRange.Replace("6.000,00", "6000.00", 1) | Out-Null
$importeRange.Replace("230,00", "230.00", 1) | Out-Null
$importeRange.Replace("400,00", "400.00", 1) | Out-Null
$importeRange.Replace("124.988,05", "124988.05", 1) | Out-Null
$importeRange.Replace("4.125,00", "4125.00", 1) | Out-Null
$importeRange.Replace("1.400,00", "1400.00", 1) | Out-Null
$importeRange.Replace("16,00", "16.00", 1) | Out-Null
$importeRange.Replace("504,50", "504.50", 1) | Out-Null
$importeRange.Replace("1.367,45", "1367.45", 1) | Out-Null
$importeRange.Replace("236.726,79", "236726.79", 1) | Out-Null
$importeRange.Replace("10.500,00", "10500.00", 1) | Out-Null
$importeRange.Replace("220,00", "220.00", 1) | Out-Null
$importeRange.Replace("1.675,00", "1675.00", 1) | Out-Null
$importeRange.Replace("650,00", "650.00", 1) | Out-Null
$importeRange.Replace("5.745,24", "5745.24", 1) | Out-Null
$importeRange.Replace("200,00", "200.00", 1) | Out-Null
$importeRange.Replace("250,00", "250.00", 1) | Out-Null
$importeRange.Replace("2.500,00", "2500.00", 1) | Out-Null
$importeRange.Replace("25.000,00", "25000.00", 1) | Out-Null
$importeRange.Replace("120,00", "120.00", 1) | Out-Null
$importeRange.Replace("7.969,35", "7969.35", 1) | Out-Null
$importeRange.Replace("85,00", "85.00", 1) | Out-Null
$importeRange.Replace("50,00", "50.00", 1) | Out-Null
$importeRange.Replace("20,00", "20.00", 1) | Out-Null
$importeRange.Replace("95,00", "95.00", 1) | Out-Null
$importeRange.Replace("125,00", "125.00", 1) | Out-Null
$importeRange.Replace("1.065,00", "1065.00", 1) | Out-Null
$importeRange.Replace("987,36", "987.36", 1) | Out-Null
$importeRange.Replace("3.798,62", "3798.62", 1) | Out-Null
$importeRange.Replace("425,00", "425.00", 1) | Out-Null
$importeRange.Replace("215,10", "215.10", 1) | Out-Null
$importeRange.Replace("736,00", "736.00", 1) | Out-Null
$importeRange.Replace("75,00", "75.00", 1) | Out-Null
$importeRange.Replace("10.002,20", "10002.20", 1) | Out-Null
$importeRange.Replace("952,38", "952.38", 1) | Out-Null
$importeRange.Replace("19.282,12", "19282.12", 1) | Out-Null
$importeRange.Replace("2.360,00", "2360.00", 1) | Out-Null
$importeRange.Replace("477,75", "477.75", 1) | Out-Null
$importeRange.Replace("2.913,00", "2913.00", 1) | Out-Null
$importeRange.Replace("8.625,00", "8625.00", 1) | Out-Null
$importeRange.Replace("729,25", "729.25", 1) | Out-Null
$importeRange.Replace("2.124,20", "2124.20", 1) | Out-Null
$importeRange.Replace("3.890,00", "3890.00", 1) | Out-Null
$importeRange.Replace("530,00", "530.00", 1) | Out-Null
$importeRange.Replace("450,00", "450.00", 1) | Out-Null
$importeRange.Style = "Normal"

Write-Host "done"